# Add "RDP" and "File transfer" entries to the Portal sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Portal")

# New row 48: RDP
$ws.Range("A48").Value = "RDP"
$ws.Range("B48").Value = "Win10 RDPs"
$ws.Range("C48").Value = "d99, d100, t66, t67, p21, p22, p26, u27"
$ws.Range("D48").Value = 3389

# New row 49: File transfer
$ws.Range("A49").Value = "File transfer"
$ws.Range("B49").Value = "Win10 RDPs"
$ws.Range("C49").Value = "d99, d100, t66, t67, p21, p22, p26, u27"
$ws.Range("D49").Value = "137-139, 445"
$ws.Range("E49").Value = "udp & tcp"

# Column D got a bit wider to fit the new content (~11.7 characters).
$ws.Columns.Item(4).ColumnWidth = 10.8

# Update the view: move the active selection to the new last row
# (matches the author's saved view state).
$ws.Activate()
$ws.Range("E50").Select()
